$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold text-formatted numbers (e.g. "158.40", "65.658.84").
# Force text format so Excel does not coerce/normalize them into real numbers
# (which would drop formatting like trailing zeros or thousand-dot separators).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.649.62'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.667.85'
$ws.Range("E3").Value = '  -0.95%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.18'
$ws.Range("E5").Value = '  -1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.49'
$ws.Range("E6").Value = '  -0.77%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +4.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.129'
$ws.Range("E9").Value = '  +4.62%  '

$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("E11").Value = '  -2.45%  '

$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.28'
$ws.Range("E13").Value = '  -3.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000196'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.147.14'
$ws.Range("E15").Value = '  -1.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.493.07'
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.671.85'
$ws.Range("E17").Value = '  -0.81%  '

$ws.Range("E18").Value = '  +1.32%  '

$ws.Range("E19").Value = '  -2.08%  '

$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.56'
$ws.Range("E21").Value = '  -2.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.49'
$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000111'
$ws.Range("E24").Value = '  +4.58%  '

$ws.Range("E25").Value = '  -2.48%  '

$ws.Range("E26").Value = '  -2.50%  '

$ws.Range("E27").Value = '  -2.62%  '

$ws.Range("E28").Value = '  -5.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.99'
$ws.Range("E29").Value = '  -5.05%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '526.13'
$ws.Range("E32").Value = '  -2.58%  '

$ws.Range("E33").Value = '  -2.44%  '

$ws.Range("E34").Value = '  -3.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.49'
$ws.Range("E35").Value = '  +1.77%  '

$ws.Range("E36").Value = '  -2.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.52'
$ws.Range("E37").Value = '  -1.54%  '

$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '158.40'
$ws.Range("E39").Value = '  -3.03%  '

$ws.Range("E40").Value = '  -3.02%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '163.31'
$ws.Range("E42").Value = '  -5.08%  '

$ws.Range("E43").Value = '  -1.46%  '

$ws.Range("E44").Value = '  +1.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0609'
$ws.Range("E45").Value = '  -0.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.73'
$ws.Range("E46").Value = '  -3.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.638'
$ws.Range("E47").Value = '  -2.71%  '

$ws.Range("E48").Value = '  -3.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0262'
$ws.Range("E49").Value = '  +13.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0999'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.07'
$ws.Range("E51").Value = '  -4.20%  '
